# Fruta / hortaliza, semanal
#
# The published market-price data in this sheet was re-sorted (the
# underlying daily/weekly records were re-ordered chronologically); every
# data row (columns A:R) that used to live at one row now lives at
# another. No cell's own content changes other than by relocating with
# its row. We snapshot every data row first, then write each snapshot
# back out to its new row, so the permutation applies atomically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# before-row -> after-row mapping for the data rows (2-21; row 1 is the header)
$rowMap = @{
    2  = 6
    3  = 18
    4  = 2
    5  = 14
    6  = 15
    7  = 12
    8  = 13
    9  = 4
    10 = 16
    11 = 17
    12 = 5
    13 = 3
    14 = 19
    15 = 9
    16 = 10
    17 = 11
    18 = 7
    19 = 8
    20 = 21
    21 = 20
}

# Snapshot every data row (A:R) before writing anything back, since several
# rows trade places with each other.
$snapshot = @{}
foreach ($srcRow in $rowMap.Keys) {
    $rng = $ws.Range("A$srcRow`:R$srcRow")
    $snapshot[$srcRow] = $rng.Value2
}

# Write each snapshot to its destination row.
foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $ws.Range("A$dstRow`:R$dstRow").Value2 = $snapshot[$srcRow]
}
